# Append the new daily "backup" rows (2-8) to the RDMARCAS list sheet.
# Rows 2-7 hold plain numeric values; the final row (8) was produced by a
# formatting pass (see commit message) that turned the numbers into
# formatted text (kept here as literal, already-formatted strings), so we
# force those cells to Text before assigning them to stop Excel from
# re-parsing the digits back into floats and dropping the trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "14/07/2023"
$ws.Range("B2").Value = 3000.6
$ws.Range("C2").Value = 3000.6
$ws.Range("D2").Value = 6000
$ws.Range("E2").Value = 6000
$ws.Range("F2").Value = 2999.4
$ws.Range("G2").Value = 199.96

# Row 3
$ws.Range("A3").Value = "14/07/2023"
$ws.Range("B3").Value = 9800
$ws.Range("C3").Value = 12800.6
$ws.Range("D3").Value = 12500
$ws.Range("E3").Value = 18500
$ws.Range("F3").Value = 5699.4
$ws.Range("G3").Value = 144.52

# Row 4
$ws.Range("A4").Value = "14/07/2023"
$ws.Range("B4").Value = 9000
$ws.Range("C4").Value = 21800.6
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 19500
$ws.Range("F4").Value = 2300.6
$ws.Range("G4").Value = 89.45

# Row 5
$ws.Range("A5").Value = "14/07/2023"
$ws.Range("B5").Value = 6000
$ws.Range("C5").Value = 27800.6
$ws.Range("D5").Value = 15000
$ws.Range("E5").Value = 34500
$ws.Range("F5").Value = 6699.4
$ws.Range("G5").Value = 124.1

# Row 6
$ws.Range("A6").Value = "14/07/2023"
$ws.Range("B6").Value = 6700
$ws.Range("C6").Value = 34500.6
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 34500
$ws.Range("F6").Value = 0.6
$ws.Range("G6").Value = 100

# Row 7
$ws.Range("A7").Value = "14/07/2023"
$ws.Range("B7").Value = 0.6
$ws.Range("C7").Value = 34501.2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 34500
$ws.Range("F7").Value = 1.2
$ws.Range("G7").Value = 100

# Row 8 - formatted ("R$" style backup) values stored as text so the
# trailing zeros / decimal formatting survive instead of collapsing back
# into plain floats.
$ws.Range("A8").Value = "14/07/2023"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "1.20"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "34502.40"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.00"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "34500.00"

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2.40"

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "99.99"
